$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1,1).Value = "JobTitle"
$ws.Cells.Item(1,2).Value = "Company"
$ws.Cells.Item(1,3).Value = "Location"
$ws.Cells.Item(1,4).Value = "JobType"
$ws.Cells.Item(1,5).Value = "CompanySize"
$ws.Cells.Item(1,6).Value = "DaysPostedAgo"

$ws.Cells.Item(2,1).Value = "Entry-Level Developer - required to work in office"
$ws.Cells.Item(2,2).Value = "Vaco Technology"
$ws.Cells.Item(2,3).Value = "Virginia Beach, VA"
$ws.Cells.Item(2,4).Value = "Full Time"
$ws.Cells.Item(2,5).Value = "N/A"
$ws.Cells.Item(2,6).Value = "16 days ago"

$ws.Cells.Item(3,1).Value = "Work From Home Entry Level Data Entry Associate"
$ws.Cells.Item(3,2).Value = "Level Up Entry"
$ws.Cells.Item(3,3).Value = "Philadelphia, PA, PA"
$ws.Cells.Item(3,4).Value = "Full-Time/Part-Time"
$ws.Cells.Item(3,5).Value = "N/A"
$ws.Cells.Item(3,6).Value = "5 days ago"

$ws.Cells.Item(4,1).Value = "Entry-level Healthcare Digital Technology Developer"
$ws.Cells.Item(4,2).Value = "Cognizant Technology"
$ws.Cells.Item(4,3).Value = "Dallas, TX"
$ws.Cells.Item(4,4).Value = "Full Time"
$ws.Cells.Item(4,5).Value = "N/A"
$ws.Cells.Item(4,6).Value = "20 days ago"

$ws.Cells.Item(5,1).Value = "Entry Level Software Developer - Dev10 Technology Development Program - NY"
$ws.Cells.Item(5,2).Value = "Genesis10"
$ws.Cells.Item(5,3).Value = "New York, NY"
$ws.Cells.Item(5,4).Value = "Full-Time/Part-Time"
$ws.Cells.Item(5,5).Value = "N/A"
$ws.Cells.Item(5,6).Value = "28 days ago"

$ws.Cells.Item(6,1).Value = "Entry Level Software Developer - Dev10 Technology Development Program - WI"
$ws.Cells.Item(6,2).Value = "Genesis10"
$ws.Cells.Item(6,3).Value = "Milwaukee, WI"
$ws.Cells.Item(6,4).Value = "Full Time"
$ws.Cells.Item(6,5).Value = "N/A"
$ws.Cells.Item(6,6).Value = "Today"

$ws.Cells.Item(7,1).Value = "Entry Level Software Developer - Dev10 Technology Development Program - MN"
$ws.Cells.Item(7,2).Value = "Genesis10"
$ws.Cells.Item(7,3).Value = "Minneapolis, MN"
$ws.Cells.Item(7,4).Value = "Full Time"
$ws.Cells.Item(7,5).Value = "N/A"
$ws.Cells.Item(7,6).Value = "7 days ago"

$ws.Cells.Item(8,1).Value = "Entry-level EAS Digital Technology Developer"
$ws.Cells.Item(8,2).Value = "Cognizant Technology"
$ws.Cells.Item(8,3).Value = "Dallas, TX"
$ws.Cells.Item(8,4).Value = "Full Time"
$ws.Cells.Item(8,5).Value = "N/A"
$ws.Cells.Item(8,6).Value = "11 days ago"

$ws.Cells.Item(9,1).Value = "Entry Level Software Developer - Dev10 Technology Development Program - .NET - TX"
$ws.Cells.Item(9,2).Value = "Genesis10"
$ws.Cells.Item(9,3).Value = "Dallas, TX"
$ws.Cells.Item(9,4).Value = "Full Time"
$ws.Cells.Item(9,5).Value = "N/A"
$ws.Cells.Item(9,6).Value = "11 days ago"

$ws.Cells.Item(10,1).Value = "Entry Level Software Developer - Dev10 Technology Development Program - .NET - MN"
$ws.Cells.Item(10,2).Value = "Genesis10"
$ws.Cells.Item(10,3).Value = "Minneapolis, MN"
$ws.Cells.Item(10,4).Value = "Full Time"
$ws.Cells.Item(10,5).Value = "N/A"
$ws.Cells.Item(10,6).Value = "11 days ago"

$ws.Cells.Item(11,1).Value = "Software Developer - Entry Level"
$ws.Cells.Item(11,2).Value = "Revature"
$ws.Cells.Item(11,3).Value = "Reston, VA"
$ws.Cells.Item(11,4).Value = "Full Time"
$ws.Cells.Item(11,5).Value = "N/A"
$ws.Cells.Item(11,6).Value = "21 days ago"

$ws.Cells.Item(12,1).Value = "Entry-Level Developer - required to work in office"
$ws.Cells.Item(12,2).Value = "Vaco Technology"
$ws.Cells.Item(12,3).Value = "Virginia Beach, VA"
$ws.Cells.Item(12,4).Value = "Full Time"
$ws.Cells.Item(12,5).Value = "N/A"
$ws.Cells.Item(12,6).Value = "TEMP"
$ws.Cells.Item(12,6).Value = ""

$ws.Cells.Item(13,1).Value = "Work From Home Entry Level Data Entry Associate"
$ws.Cells.Item(13,2).Value = "Level Up Entry"
$ws.Cells.Item(13,3).Value = "Philadelphia, PA, PA"
$ws.Cells.Item(13,4).Value = "Full-Time/Part-Time"
$ws.Cells.Item(13,5).Value = "N/A"
$ws.Cells.Item(13,6).Value = "TEMP"
$ws.Cells.Item(13,6).Value = ""

$ws.Cells.Item(14,1).Value = "Entry-Level Developer - required to work in office"
$ws.Cells.Item(14,2).Value = "Vaco Technology"
$ws.Cells.Item(14,3).Value = "Virginia Beach, VA"
$ws.Cells.Item(14,4).Value = "Full Time"
$ws.Cells.Item(14,5).Value = "N/A"
$ws.Cells.Item(14,6).Value = "TEMP"
$ws.Cells.Item(14,6).Value = ""

$ws.Cells.Item(15,1).Value = "Entry-Level Developer - required to work in office"
$ws.Cells.Item(15,2).Value = "Vaco Technology"
$ws.Cells.Item(15,3).Value = "Virginia Beach, VA"
$ws.Cells.Item(15,4).Value = "Full Time"
$ws.Cells.Item(15,5).Value = "N/A"
$ws.Cells.Item(15,6).Value = "TEMP"
$ws.Cells.Item(15,6).Value = ""

$ws.Cells.Item(16,1).Value = "Work From Home Entry Level Data Entry Associate"
$ws.Cells.Item(16,2).Value = "Level Up Entry"
$ws.Cells.Item(16,3).Value = "Philadelphia, PA, PA"
$ws.Cells.Item(16,4).Value = "Full-Time/Part-Time"
$ws.Cells.Item(16,5).Value = "N/A"
$ws.Cells.Item(16,6).Value = "TEMP"
$ws.Cells.Item(16,6).Value = ""

$ws.Cells.Item(17,1).Value = "Entry-level Healthcare Digital Technology Developer"
$ws.Cells.Item(17,2).Value = "Cognizant Technology"
$ws.Cells.Item(17,3).Value = "Dallas, TX"
$ws.Cells.Item(17,4).Value = "Full Time"
$ws.Cells.Item(17,5).Value = "N/A"
$ws.Cells.Item(17,6).Value = "TEMP"
$ws.Cells.Item(17,6).Value = ""

$ws.Cells.Item(18,1).Value = "Entry-Level Developer - required to work in office"
$ws.Cells.Item(18,2).Value = "Vaco Technology"
$ws.Cells.Item(18,3).Value = "Virginia Beach, VA"
$ws.Cells.Item(18,4).Value = "Full Time"
$ws.Cells.Item(18,5).Value = "N/A"
$ws.Cells.Item(18,6).Value = "TEMP"
$ws.Cells.Item(18,6).Value = ""

$ws.Cells.Item(19,1).Value = "Work From Home Entry Level Data Entry Associate"
$ws.Cells.Item(19,2).Value = "Level Up Entry"
$ws.Cells.Item(19,3).Value = "Philadelphia, PA, PA"
$ws.Cells.Item(19,4).Value = "Full-Time/Part-Time"
$ws.Cells.Item(19,5).Value = "N/A"
$ws.Cells.Item(19,6).Value = "TEMP"
$ws.Cells.Item(19,6).Value = ""

$ws.Cells.Item(20,1).Value = "Entry-level Healthcare Digital Technology Developer"
$ws.Cells.Item(20,2).Value = "Cognizant Technology"
$ws.Cells.Item(20,3).Value = "Dallas, TX"
$ws.Cells.Item(20,4).Value = "Full Time"
$ws.Cells.Item(20,5).Value = "N/A"
$ws.Cells.Item(20,6).Value = "TEMP"
$ws.Cells.Item(20,6).Value = ""

$ws.Cells.Item(21,1).Value = "Entry-Level Developer - required to work in office"
$ws.Cells.Item(21,2).Value = "Vaco Technology"
$ws.Cells.Item(21,3).Value = "Virginia Beach, VA"
$ws.Cells.Item(21,4).Value = "Full Time"
$ws.Cells.Item(21,5).Value = "N/A"
$ws.Cells.Item(21,6).Value = "TEMP"
$ws.Cells.Item(21,6).Value = ""

$ws.Cells.Item(22,1).Value = "Work From Home Entry Level Data Entry Associate"
$ws.Cells.Item(22,2).Value = "Level Up Entry"
$ws.Cells.Item(22,3).Value = "Philadelphia, PA, PA"
$ws.Cells.Item(22,4).Value = "Full-Time/Part-Time"
$ws.Cells.Item(22,5).Value = "N/A"
$ws.Cells.Item(22,6).Value = "TEMP"
$ws.Cells.Item(22,6).Value = ""

$ws.Cells.Item(23,1).Value = "Entry-level Healthcare Digital Technology Developer"
$ws.Cells.Item(23,2).Value = "Cognizant Technology"
$ws.Cells.Item(23,3).Value = "Dallas, TX"
$ws.Cells.Item(23,4).Value = "Full Time"
$ws.Cells.Item(23,5).Value = "N/A"
$ws.Cells.Item(23,6).Value = "TEMP"
$ws.Cells.Item(23,6).Value = ""

$ws.Cells.Item(24,1).Value = "Entry-Level Developer - required to work in office"
$ws.Cells.Item(24,2).Value = "Vaco Technology"
$ws.Cells.Item(24,3).Value = "Virginia Beach, VA"
$ws.Cells.Item(24,4).Value = "Full Time"
$ws.Cells.Item(24,5).Value = "N/A"
$ws.Cells.Item(24,6).Value = "TEMP"
$ws.Cells.Item(24,6).Value = ""

$ws.Cells.Item(25,1).Value = "Work From Home Entry Level Data Entry Associate"
$ws.Cells.Item(25,2).Value = "Level Up Entry"
$ws.Cells.Item(25,3).Value = "Philadelphia, PA, PA"
$ws.Cells.Item(25,4).Value = "Full-Time/Part-Time"
$ws.Cells.Item(25,5).Value = "N/A"
$ws.Cells.Item(25,6).Value = "TEMP"
$ws.Cells.Item(25,6).Value = ""

$ws.Cells.Item(26,1).Value = "Entry-level Healthcare Digital Technology Developer"
$ws.Cells.Item(26,2).Value = "Cognizant Technology"
$ws.Cells.Item(26,3).Value = "Dallas, TX"
$ws.Cells.Item(26,4).Value = "Full Time"
$ws.Cells.Item(26,5).Value = "N/A"
$ws.Cells.Item(26,6).Value = "TEMP"
$ws.Cells.Item(26,6).Value = ""

$ws.Cells.Item(27,1).Value = "Entry Level Software Developer - Dev10 Technology Development Program - NY"
$ws.Cells.Item(27,2).Value = "Genesis10"
$ws.Cells.Item(27,3).Value = "New York, NY"
$ws.Cells.Item(27,4).Value = "Full-Time/Part-Time"
$ws.Cells.Item(27,5).Value = "N/A"
$ws.Cells.Item(27,6).Value = "TEMP"
$ws.Cells.Item(27,6).Value = ""

$ws.Cells.Item(28,1).Value = "Entry Level Software Developer - Dev10 Technology Development Program - WI"
$ws.Cells.Item(28,2).Value = "Genesis10"
$ws.Cells.Item(28,3).Value = "Milwaukee, WI"
$ws.Cells.Item(28,4).Value = "Full Time"
$ws.Cells.Item(28,5).Value = "N/A"
$ws.Cells.Item(28,6).Value = "TEMP"
$ws.Cells.Item(28,6).Value = ""

$ws.Cells.Item(29,1).Value = "Entry Level Software Developer - Dev10 Technology Development Program - MN"
$ws.Cells.Item(29,2).Value = "Genesis10"
$ws.Cells.Item(29,3).Value = "Minneapolis, MN"
$ws.Cells.Item(29,4).Value = "Full Time"
$ws.Cells.Item(29,5).Value = "N/A"
$ws.Cells.Item(29,6).Value = "TEMP"
$ws.Cells.Item(29,6).Value = ""

$ws.Cells.Item(30,1).Value = "Entry-level EAS Digital Technology Developer"
$ws.Cells.Item(30,2).Value = "Cognizant Technology"
$ws.Cells.Item(30,3).Value = "Dallas, TX"
$ws.Cells.Item(30,4).Value = "Full Time"
$ws.Cells.Item(30,5).Value = "N/A"
$ws.Cells.Item(30,6).Value = "TEMP"
$ws.Cells.Item(30,6).Value = ""

$ws.Cells.Item(31,1).Value = "Entry Level Software Developer - Dev10 Technology Development Program - .NET - TX"
$ws.Cells.Item(31,2).Value = "Genesis10"
$ws.Cells.Item(31,3).Value = "Dallas, TX"
$ws.Cells.Item(31,4).Value = "Full Time"
$ws.Cells.Item(31,5).Value = "N/A"
$ws.Cells.Item(31,6).Value = "TEMP"
$ws.Cells.Item(31,6).Value = ""

$ws.Cells.Item(32,1).Value = "Entry Level Software Developer - Dev10 Technology Development Program - .NET - MN"
$ws.Cells.Item(32,2).Value = "Genesis10"
$ws.Cells.Item(32,3).Value = "Minneapolis, MN"
$ws.Cells.Item(32,4).Value = "Full Time"
$ws.Cells.Item(32,5).Value = "N/A"
$ws.Cells.Item(32,6).Value = "TEMP"
$ws.Cells.Item(32,6).Value = ""

$ws.Cells.Item(33,1).Value = "Software Developer - Entry Level"
$ws.Cells.Item(33,2).Value = "Revature"
$ws.Cells.Item(33,3).Value = "Reston, VA"
$ws.Cells.Item(33,4).Value = "Full Time"
$ws.Cells.Item(33,5).Value = "N/A"
$ws.Cells.Item(33,6).Value = "TEMP"
$ws.Cells.Item(33,6).Value = ""

$ws.Cells.Item(34,1).Value = "Entry-Level Developer - required to work in office"
$ws.Cells.Item(34,2).Value = "Vaco Technology"
$ws.Cells.Item(34,3).Value = "Virginia Beach, VA"
$ws.Cells.Item(34,4).Value = "Full Time"
$ws.Cells.Item(34,5).Value = "N/A"
$ws.Cells.Item(34,6).Value = "TEMP"
$ws.Cells.Item(34,6).Value = ""

$ws.Cells.Item(35,1).Value = "Entry-level Healthcare Digital Technology Developer"
$ws.Cells.Item(35,2).Value = "Cognizant Technology"
$ws.Cells.Item(35,3).Value = "Dallas, TX"
$ws.Cells.Item(35,4).Value = "Full Time"
$ws.Cells.Item(35,5).Value = "N/A"
$ws.Cells.Item(35,6).Value = "TEMP"
$ws.Cells.Item(35,6).Value = ""

$ws.Cells.Item(36,1).Value = "Entry-level EAS Digital Technology Developer"
$ws.Cells.Item(36,2).Value = "Cognizant Technology"
$ws.Cells.Item(36,3).Value = "Dallas, TX"
$ws.Cells.Item(36,4).Value = "Full Time"
$ws.Cells.Item(36,5).Value = "N/A"
$ws.Cells.Item(36,6).Value = "TEMP"
$ws.Cells.Item(36,6).Value = ""

$ws.Cells.Item(37,1).Value = "Entry Level SQL Developer"
$ws.Cells.Item(37,2).Value = "HAWAII MAINLAND ADMINISTRATORS L"
$ws.Cells.Item(37,3).Value = "Tempe, AZ"
$ws.Cells.Item(37,4).Value = "Full Time"
$ws.Cells.Item(37,5).Value = "N/A"
$ws.Cells.Item(37,6).Value = "TEMP"
$ws.Cells.Item(37,6).Value = ""

$ws.Cells.Item(38,1).Value = "Entry-Level .NET Developer / Application Support"
$ws.Cells.Item(38,2).Value = "Medline Industries, Inc."
$ws.Cells.Item(38,3).Value = "Mundelein, IL"
$ws.Cells.Item(38,4).Value = "Full Time"
$ws.Cells.Item(38,5).Value = "N/A"
$ws.Cells.Item(38,6).Value = "TEMP"
$ws.Cells.Item(38,6).Value = ""

$ws.Cells.Item(39,1).Value = "Entry Level Unreal Engine C++ Developer"
$ws.Cells.Item(39,2).Value = "Opex"
$ws.Cells.Item(39,3).Value = "Moorestown, NJ"
$ws.Cells.Item(39,4).Value = "Full Time"
$ws.Cells.Item(39,5).Value = "N/A"
$ws.Cells.Item(39,6).Value = "TEMP"
$ws.Cells.Item(39,6).Value = ""

$ws.Cells.Item(40,1).Value = "Entry-Level Developer - required to work in office"
$ws.Cells.Item(40,2).Value = "Vaco Technology"
$ws.Cells.Item(40,3).Value = "Virginia Beach, VA"
$ws.Cells.Item(40,4).Value = "Full Time"
$ws.Cells.Item(40,5).Value = "N/A"
$ws.Cells.Item(40,6).Value = "TEMP"
$ws.Cells.Item(40,6).Value = ""

$ws.Cells.Item(41,1).Value = "Work From Home Entry Level Data Entry Associate"
$ws.Cells.Item(41,2).Value = "Level Up Entry"
$ws.Cells.Item(41,3).Value = "Philadelphia, PA, PA"
$ws.Cells.Item(41,4).Value = "Full-Time/Part-Time"
$ws.Cells.Item(41,5).Value = "N/A"
$ws.Cells.Item(41,6).Value = "TEMP"
$ws.Cells.Item(41,6).Value = ""

$ws.Cells.Item(42,1).Value = "Entry-level Healthcare Digital Technology Developer"
$ws.Cells.Item(42,2).Value = "Cognizant Technology"
$ws.Cells.Item(42,3).Value = "Dallas, TX"
$ws.Cells.Item(42,4).Value = "Full Time"
$ws.Cells.Item(42,5).Value = "N/A"
$ws.Cells.Item(42,6).Value = "TEMP"
$ws.Cells.Item(42,6).Value = ""

$ws.Cells.Item(43,1).Value = "Entry Level Marketing"
$ws.Cells.Item(43,2).Value = "DFW Brands"
$ws.Cells.Item(43,3).Value = "Dallas, TX"
$ws.Cells.Item(43,4).Value = "Full Time"
$ws.Cells.Item(43,5).Value = "N/A"
$ws.Cells.Item(43,6).Value = "TEMP"
$ws.Cells.Item(43,6).Value = ""

$ws.Cells.Item(44,1).Value = "Entry Level Engineer"
$ws.Cells.Item(44,2).Value = "Prokatchers LLC"
$ws.Cells.Item(44,3).Value = "Irving, TX"
$ws.Cells.Item(44,4).Value = "Contractor"
$ws.Cells.Item(44,5).Value = "N/A"
$ws.Cells.Item(44,6).Value = "TEMP"
$ws.Cells.Item(44,6).Value = ""

$ws.Cells.Item(45,1).Value = "Facilities Engineer - Electrical - Entry/Experienced Level (NSAW and NSAH)"
$ws.Cells.Item(45,2).Value = "National Security Agency (NSA)"
$ws.Cells.Item(45,3).Value = "Fort Meade, MD"
$ws.Cells.Item(45,4).Value = "Full Time"
$ws.Cells.Item(45,5).Value = "N/A"
